$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16 (this shifts the old row 16 -> 17, old row 17 -> 18,
# and shifts the merged cell refs on those rows automatically).
$ws.Rows.Item(16).Insert()

# Populate the new product row with the new sale line's data.
$ws.Range("A16").Value = 10
$ws.Range("C16").Value = "مناديل سولو سحب"
$ws.Range("H16").Value = "20:0"
# L16 and P16 hold numeric-looking text ("0" and "45.0000") that must be stored as
# text, not get auto-converted to numbers, so enter them with a leading apostrophe.
$ws.Range("L16").Value = "'0"
$ws.Range("N16").Value = "45.00"
$ws.Range("P16").Value = "'45.0000"
$ws.Range("Q16").Value = "1:0"

# Copy formatting (styles/borders/fills/fonts/number formats) from row 15 (an existing
# product row) onto the new row 16 so it matches the look of the other data rows -
# this also clears the quote-prefix marker left behind by the apostrophe entries above.
$ws.Range("A15:Q15").Copy()
$ws.Range("A16:Q16").PasteSpecial(-4122)
$ws.Rows.Item(16).RowHeight = 25.5

# Re-create the merged cells for the new product row, matching the pattern used by
# every other product row (A:B, C:G, H:K, L:M, N:O merged; P and Q stay separate).
$ws.Range("A16:B16").Merge()
$ws.Range("C16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()
$ws.Range("N16:O16").Merge()

# Update the running total (old row16/P16 shifted down to row17/P17): 772 -> 817.
$ws.Range("P17").Value = 817

# Update the footer timestamp (old row17/A17 shifted down to row18/A18) to reflect the
# new generation time.
$ws.Range("A18").Value = "Monday, 22 September, 2025 10:37 AM"
